$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53 (pushes existing rows 53..158 down to 54..159)
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new record
$ws.Range("A53").Value = 11
$ws.Range("B53").Value = "Vega Monumental Concepción"
$ws.Range("C53").Value = "Bíobío"
$ws.Range("D53").Value = 45272
$ws.Range("E53").Value = 8
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100101
$ws.Range("H53").Value = "Berries"
$ws.Range("I53").Value = 100101001
$ws.Range("J53").Value = "Arándano (blue)"
$ws.Range("K53").Value = "Sin especificar"
$ws.Range("L53").Value = "Primera"
$ws.Range("M53").Value = 100
$ws.Range("N53").Value = 4000
$ws.Range("O53").Value = 4500
$ws.Range("P53").Value = 4250
$ws.Range("Q53").Value = "$/bandeja 2 kilos"
$ws.Range("R53").Value = "Región de Ñuble"
$ws.Range("S53").Value = 2125
$ws.Range("T53").Value = 2
